# Improve: Copy Module hint
# Insert a new "MailTemplatePathCopyModule" settings row right above the
# existing "ExcelOfferThreshhold" row (new row 18) on the "Settings" sheet,
# shifting the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new, blank row at row 18 (everything from row 18 down shifts to 19+).
$ws.Rows.Item(18).Insert()

# Populate the new row with the Copy-Module mail-template hint, reusing the
# same "relative path" helper text already used by the sibling rows above.
$ws.Range("A18").Value = "MailTemplatePathCopyModule"
$ws.Range("B18").Value = "Data\Input\MailTemplateCopyModule"
$ws.Range("C18").Value = $ws.Range("C17").Value2

# Match the row height used by its neighbouring "path" rows (17/19).
$ws.Rows.Item(18).RowHeight = 30

# Leave the new cell selected, as in the authored workbook.
[void]$ws.Range("C18").Select()
